# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes    = $wb.Worksheets.Item("全部类型")

# Row -> new F value, for the "展览" sheet
$exhibitionUpdates = @{
    2  = 7056
    7  = 155
    8  = 120
    12 = 200
    15 = 1829
    17 = 3653
    22 = 28
    23 = 2280
    27 = 35
    31 = 158
    32 = 1294
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Row -> new F value, for the "全部类型" sheet
$allTypesUpdates = @{
    2  = 7056
    8  = 155
    9  = 120
    13 = 200
    16 = 1829
    18 = 3653
    23 = 28
    24 = 2280
    28 = 35
    32 = 158
    33 = 1294
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
